# Append the 2023/24 wild report harvest rows (rows 2-10) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows use a plain (non-bordered) style with wrap-text turned on,
# matching the added "s=2" cellXfs entry.
$dataRange = $ws.Range("A2:M10")
$dataRange.WrapText = $true

# Columns: A Trapping Licence Year, B Did Harvest Occur?, C Trapline Type,
# D Trapline Number, E Month, F Species, G WMU, H Male Count, I Female Count,
# J Unknown Sex Count, K Harvest in Park?, L Park Name, M PERMITAUTHORIZATIONNUMBER
$rows = @(
    @("2023/24","Yes","Registered Trapline","TR0331T005","November","Squirrel","331",0,0,5,"No","",""),
    @("2023/24","Yes","Registered Trapline","TR0331T005","November","Flying Squirrel","331",0,0,3,"No","",""),
    @("2023/24","Yes","Registered Trapline","TR0331T005","November","Squirrel","331",0,0,2,"Yes","Marble Range Park","FILL IN WITH PERMIT AUTHORIZATION NUMBER"),
    @("2023/24","Yes","Registered Trapline","TR0331T005","December","Weasel","331",0,0,1,"No","",""),
    @("2023/24","Yes","Registered Trapline","TR0331T005","December","Fisher","331",1,0,0,"No","",""),
    @("2023/24","Yes","Registered Trapline","TR0331T005","January","Squirrel","331",0,0,15,"No","",""),
    @("2023/24","Yes","Registered Trapline","TR0331T005","January","Weasel","331",0,0,2,"No","",""),
    @("2023/24","Yes","Registered Trapline","TR0331T005","January","Marten","331",1,0,0,"No","",""),
    @("2023/24","Yes","Registered Trapline","TR0331T005","January","Squirrel","331",0,0,1,"Yes","Marble Range Park","FILL IN WITH PERMIT AUTHORIZATION NUMBER")
)

$numericCols = @(8,9,10)   # H, I, J are numeric Male/Female/Unknown counts

$r = 2
foreach ($row in $rows) {
    for ($c = 1; $c -le $row.Length; $c++) {
        $val = $row[$c - 1]
        if ($numericCols -contains $c) {
            $ws.Cells.Item($r, $c).Value = $val
        } else {
            $ws.Cells.Item($r, $c).Value = [string]$val
        }
    }
    $r++
}
